$d = $word.ActiveDocument

function Set-ParagraphText($para, $text) {
    # Replace the paragraph's run content (excluding the trailing paragraph
    # mark) with a single merged run containing $text. Clearing first avoids
    # a same-content-so-no-op short circuit, ensuring the runs actually
    # collapse into one.
    $r = $para.Range
    $r.End = $r.End - 1
    $r.Text = ""
    $r2 = $para.Range
    $r2.End = $r2.End - 1
    $r2.Text = $text
}

# Title paragraph: "Answers:" " " "Introduction" " " "to" " " "radians"
#   -> single run "Answers: Introduction to radians"
Set-ParagraphText $d.Paragraphs.Item(1) "Answers: Introduction to radians"

# Author paragraph: "Ifan" " " "Howells-Baines," " " "Mark" " " "Toner"
#   -> single run "Ifan Howells-Baines, Mark Toner"
Set-ParagraphText $d.Paragraphs.Item(2) "Ifan Howells-Baines, Mark Toner"

# Abstract paragraph: "Answers" " " "to" " " "the" ... "radians."
#   -> single run "Answers to the questions relating to the guide on radians."
Set-ParagraphText $d.Paragraphs.Item(4) "Answers to the questions relating to the guide on radians."
